$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "pH 9.2"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("pH 9.2")

# pH dilution series: last point renamed from 9 to 8.7
$ws1.Range("A13").Value = 8.7
$ws1.Range("A14").Value = 8.7

# Bicarbonate table: add "HCO3 [mg/L]" header and measured values,
# collapse the Feed/Permeate pairs down to just the Feed rows with data
$ws1.Range("D17").Value = "HCO3 [mg/L]"

$ws1.Range("D18").Value = 465

$ws1.Range("A19").Value = 6
$ws1.Range("B19").Value = "Feed"
$ws1.Range("D19").Value = 710

$ws1.Range("A20").Value = 8.7
$ws1.Range("D20").Value = 963

$ws1.Rows("21:23").Delete()

# ---------------------------------------------------------------------
# Sheet "pH 10"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pH 10")

$ws2.Range("D17").Value = "HCO3 [mg/L]"

$ws2.Rows("19:19").Delete()
$ws2.Rows("20:20").Delete()
$ws2.Rows("20:20").Delete()

# ---------------------------------------------------------------------
# Sheet "pH 10.5"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("pH 10.5")

$ws3.Range("D19").Value = "HCO3 [mg/L]"

$ws3.Rows("21:21").Delete()
$ws3.Rows("22:22").Delete()
$ws3.Rows("22:22").Delete()

# ---------------------------------------------------------------------
# Window / view bookkeeping
# ---------------------------------------------------------------------
$ws1.Select()
